$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are written as exact text (preserve trailing
# zeros / formatting) instead of being coerced to numeric values.
# NOTE: multi-area (comma-joined) ranges only apply NumberFormat to the
# first area in this runtime, so each contiguous block is set separately.
$ws.Range("D2:D19").NumberFormat = "@"
$ws.Range("D21:D26").NumberFormat = "@"
$ws.Range("D41:D45").NumberFormat = "@"
$ws.Range("D47:D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "248.00"

# Row 3
$ws.Range("D3").Value = "22.57"

# Row 4
$ws.Range("D4").Value = "5.254"

# Row 5
$ws.Range("D5").Value = "0.05685"

# Row 6
$ws.Range("D6").Value = "3.417"

# Row 7
$ws.Range("D7").Value = "6.314"

# Row 8
$ws.Range("D8").Value = "0.8066"

# Row 9
$ws.Range("D9").Value = "0.8996"

# Row 10
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01093"
$ws.Range("E10").Value = "9OneONE"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1419"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07441"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03056"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03092"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09388"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.880"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001591"
$ws.Range("E17").Value = "16BitForexTokenBF"

# Row 18
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04768"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19
$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").Value = "0.01826"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"

# Row 21
$ws.Range("D21").Value = "0.005042"

# Row 22
$ws.Range("D22").Value = "0.0009965"

# Row 23
$ws.Range("D23").Value = "0.0001500"

# Row 24
$ws.Range("D24").Value = "3.685"

# Row 25
$ws.Range("D25").Value = "2.160"

# Row 26
$ws.Range("D26").Value = "0.3245"

# Row 41
$ws.Range("D41").Value = "0.006812"

# Row 42
$ws.Range("D42").Value = "0.1064"

# Row 43
$ws.Range("D43").Value = "0.003201"

# Row 44
$ws.Range("D44").Value = "0.008753"

# Row 45
$ws.Range("D45").Value = "0.00005582"

# Row 47
$ws.Range("D47").Value = "0.4991"

# Row 48
$ws.Range("D48").Value = "0.1333"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# Row 50
$ws.Range("D50").Value = "0.01010"

